$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 ("Beton" / "P_Beton" / 9) is removed entirely; rows below shift up.
$ws.Rows.Item(7).Delete()

# The "Montageart" options string no longer references the now-removed
# P_Beton variable; it is hard-coded to the value 9 instead.
$ws.Range("D5").Value = "Betonieren:9, Konsole:1"

# Reflect the row that is now selected after the deletion (whole row 7).
$ws.Rows.Item(7).Select()
